$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.608.19"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "1.741.75"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'246.53"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "'0.4908"
$ws.Range("E7").Value = "  +2.47%  "
$ws.Range("D8").Value = "'0.2675"
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").Value = "'0.06288"
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("D10").Value = "1.748.10"
$ws.Range("D11").Value = "'0.07045"
$ws.Range("E11").Value = "  -1.10%  "
$ws.Range("D12").Value = "'15.74"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").Value = "'0.6149"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").Value = "'4.587"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").Value = "'78.09"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "26.625.77"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").Value = "'0.000007310"
$ws.Range("E18").Value = "  +5.46%  "
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "'11.58"
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("D21").Value = "1.974.60"
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("D22").Value = "'4.576"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").Value = "'8.716"
$ws.Range("E23").Value = "  -2.20%  "
$ws.Range("D24").Value = "'5.279"
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("D25").Value = "'139.31"
$ws.Range("E25").Value = "  +2.10%  "
$ws.Range("D26").Value = "'15.44"
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").Value = "'1.425"
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("D28").Value = "'1.762"
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("D29").Value = "'107.59"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("D30").Value = "'4.039"
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("D31").Value = "'0.08054"
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("D32").Value = "'3.731"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").Value = "'0.04615"
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("B34").Value = "Frax"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D34").Value = "'0.9997"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.614"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.016"
$ws.Range("E36").Value = "  +2.72%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.6405"
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'2.076"
$ws.Range("E38").Value = "  +2.91%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'0.9049"
$ws.Range("E39").Value = "  -3.27%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.427"
$ws.Range("E40").Value = "  +0.62%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.003"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.01505"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'102.11"
$ws.Range("E43").Value = "  -4.66%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.429"
$ws.Range("E44").Value = "  -3.60%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.3921"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "'6.869"
$ws.Range("E46").Value = "  -2.01%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.1185"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.05395"
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("D49").Value = "'30.63"
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.802"
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.261"
$ws.Range("E51").Value = "  -0.49%  "
